$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths: both columns become 15.42578125 chars wide (closest achievable
# quantization in this engine is 15.5, since ColumnWidth is rounded to the nearest
# 1/6th of a character as Excel does internally for pixel-based fonts).
$ws.Columns.Item(1).ColumnWidth = 14.583333333333334
$ws.Columns.Item(2).ColumnWidth = 14.583333333333334

# Cell values for A1:B32
$ws.Range("A1").Value = -0.090022772380280003
$ws.Range("B1").Value = 0.089515123139200625
$ws.Range("A2").Value = -0.067412495710039266
$ws.Range("B2").Value = 0.065276653944239094
$ws.Range("A3").Value = -0.038178160805257377
$ws.Range("B3").Value = 0.037790340319981297
$ws.Range("A4").Value = -0.029790340394653114
$ws.Range("B4").Value = 0.029444358003174464
$ws.Range("A5").Value = -0.026444358036353144
$ws.Range("B5").Value = 0.025275606112650451
$ws.Range("A6").Value = -0.0072071128811366236
$ws.Range("B6").Value = 0.0070731017464211732
$ws.Range("A7").Value = 0.0029268981609207323
$ws.Range("B7").Value = -0.0029475553045026359
$ws.Range("A8").Value = 0.012947555212585726
$ws.Range("B8").Value = -0.012969084038620515
$ws.Range("A9").Value = 0.014969084017854684
$ws.Range("B9").Value = -0.014985368582599801
$ws.Range("A10").Value = 0.016985368565189063
$ws.Range("B10").Value = -0.016985081536425639
$ws.Range("A11").Value = 0.019985081510720981
$ws.Range("B11").Value = -0.019987389303637748
$ws.Range("A12").Value = 0.023487389274704729
$ws.Range("B12").Value = -0.023532951276556524
$ws.Range("A13").Value = 0.027032951254192028
$ws.Range("B13").Value = -0.027074088015128339
$ws.Range("A14").Value = -0.0080467506084742624
$ws.Range("B14").Value = 0.0080443263646250429
$ws.Range("A15").Value = -0.0080514206072859196
$ws.Range("B15").Value = 0.0080335760280014767
$ws.Range("A16").Value = -0.0060335760357528301
$ws.Range("B16").Value = 0.0060032019058380293
$ws.Range("A17").Value = -0.0040032019142728359
$ws.Range("B17").Value = 0.003999999974335644
$ws.Range("A18").Value = -0.016102627481025422
$ws.Range("B18").Value = 0.01609118158052425
$ws.Range("A19").Value = -0.012091181615581093
$ws.Range("B19").Value = 0.012016204786439388
$ws.Range("A20").Value = -0.0080162048242637951
$ws.Range("B20").Value = 0.0080056328814901434
$ws.Range("A21").Value = -0.0040056329197684093
$ws.Range("B21").Value = 0.00399999996136291
$ws.Range("A22").Value = -0.045718770921492791
$ws.Range("B22").Value = 0.045503712374685534
$ws.Range("A23").Value = -0.040503712425364213
$ws.Range("B23").Value = 0.040099858889362849
$ws.Range("A24").Value = -0.020099859073099857
$ws.Range("B24").Value = 0.019999999813501201
$ws.Range("A25").Value = -0.016068493257209937
$ws.Range("B25").Value = 0.016010214126213285
$ws.Range("A26").Value = -0.013510214156370992
$ws.Range("B26").Value = 0.013437827808308711
$ws.Range("A27").Value = -0.010937827838752856
$ws.Range("B27").Value = 0.010524062222947439
$ws.Range("A28").Value = -0.0085240622502222863
$ws.Range("B28").Value = 0.008259801209953288
$ws.Range("A29").Value = -0.0012598012802262915
$ws.Range("B29").Value = 0.0011932048225959946
$ws.Range("A30").Value = -0.021170121357429128
$ws.Range("B30").Value = 0.021022599387158891
$ws.Range("A31").Value = -0.014022599460368212
$ws.Range("B31").Value = 0.01400099061079807
$ws.Range("A32").Value = -0.0040009907099136655
$ws.Range("B32").Value = 0.0039999999518194329
